# 2nd version - CRUDE operations for expenses and revenues and creating year
#
# Novembro: fix a bad date entry, remove a duplicated "Total gasto" box.
# Dezembro: fix a bad date entry, correct an amount, add two new expense
#           rows, switch the AMOUNT column to a plain accounting format,
#           and move the selection to where the user left off.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Novembro"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Novembro")

# "Data início:" / "Data fim:" values now show a full date-time stamp
# instead of a short date.
$ws1.Range("H1:H2").NumberFormat = "yyyy-mm-dd h:mm:ss"

# A4 held a real date (25 Oct 2023); it is corrected by hand to a plain text
# date string. "26/10/2023" isn't a valid month/day under the cell's current
# date format, so it is kept as literal text instead of being re-parsed.
$ws1.Range("A4").Value = "26/10/2023"

# The "Cadeira escritório" expense and the "Montante inicial" readout get
# the euro accounting format used elsewhere for highlighted totals.
$ws1.Range("B4").NumberFormat = "#,##0.00 €; [Red]-#,##0.00 €"
$ws1.Range("H4").NumberFormat = "#,##0.00 €; [Red]-#,##0.00 €"

# Remove the duplicated "Total gasto" box (label + SUM formula) that lived
# at G8:H8 - the sheet only keeps the original one at G7:H7.
$ws1.Range("G8:H8").Clear()

# The stray empty styled cell at I4 is removed too.
$ws1.Range("I4").Clear()

# ---------------------------------------------------------------------------
# Sheet "Dezembro"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Dezembro")

# The AMOUNT column switches from the custom euro format to a plain
# accounting format (applies to the header + the existing rows; new rows
# below inherit it automatically).
$ws2.Columns.Item(2).NumberFormat = "#,##0.00;[Red]-#,##0.00"

# A2 held a real date (25 Nov 2023); corrected by hand to plain text -
# "25/11/2023" fails to re-parse as a date (25 isn't a valid month), so it
# is kept as literal text, same as Novembro!A4 above.
$ws2.Range("A2").Value = "25/11/2023"
$ws2.Range("B2").Value = 7.3

# "Data início:" / "Data fim:" readouts lose their explicit date styling
# (they are plain text labels, not real dates).
$ws2.Range("H1").Style = "Normal"
$ws2.Range("H2").Style = "Normal"

# New expense row: Jardim zoológico / passe, 02/01/2024, 1.65
$ws2.Range("A4").NumberFormat = "@"
$ws2.Range("A4").Value = "02/01/2024"
$ws2.Range("B4").Value = 1.65
$ws2.Range("C4").Value = "Jardim zologico"
$ws2.Range("D4").Value = "passe"

# New expense row: fgsdfs / sxvvb, 04/12/2023, -21
$ws2.Range("A5").NumberFormat = "@"
$ws2.Range("A5").Value = "04/12/2023"
$ws2.Range("B5").NumberFormat = "#,##0.00 €; [Red]-#,##0.00 €"
$ws2.Range("B5").Value = -21
$ws2.Range("C5").Value = "fgsdfs"
$ws2.Range("D5").Value = "sxvvb"

# Update the selection to reflect where the user ended up working.
$ws2.Range("C12").Select()
